# Add three new ftest rows (fm32, fm33, fm34) to the "ftests" sheet,
# covering calcrules 2, 17 and 18 for excess policy layers with shares
# and a blanket deductible.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ftests")

# --- Seed formatting for the new rows by copying it from existing rows
#     whose cell styles match what we need, then overwrite the values. ---

# Row 37 needs styles: B,C,D,H,I = 5 ; E = 8 ; F,G = 7  -> matches row 18
$ws.Range("B18:I18").Copy()
$ws.Range("B37:I37").PasteSpecial(-4122)

# Rows 38/39 need styles: B,C,D,H,I = 5 ; E,F,G = 7 -> matches row 11
$ws.Range("B11:I11").Copy()
$ws.Range("B38:I38").PasteSpecial(-4122)
$ws.Range("B39:I39").PasteSpecial(-4122)

# --- Row 37: fm32 ---
$ws.Range("B37").Value = "fm32"
$ws.Range("C37").Value = "Excess policies with shares and blanket deductible (amount) using calcrule 2 only"

# --- Row 38/39 labels ---
$ws.Range("B38").Value = "fm33"
$ws.Range("B39").Value = "fm34"

# --- Row 38/39 descriptions ---
$ws.Range("C38").Value = "Excess policies with shares and blanket deductible (% loss) using calcrule 17"
$ws.Range("C39").Value = "Excess policies with shares and blanket deductible (% tiv) using calcrule 18"

# --- Row 37 remaining columns ---
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 2
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = 1
$ws.Range("H37").Value = "complete"
$ws.Range("I37").Value = "complete"

# --- Row 38 remaining columns ---
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 17
$ws.Range("F38").Value = 1
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = "complete"
$ws.Range("I38").Value = "complete"

# --- Row 39 remaining columns ---
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = 18
$ws.Range("F39").Value = 1
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = "in progress"
$ws.Range("I39").Value = "in progress"

# --- Update the sheet selection to match the new last row ---
$ws.Range("B39").Select()
